# Applies the 11:19:35 schedule-scrape refresh to Línea 141 workbook
# (sheets LP1912, LP1912-215, 6203-6173): updates header metadata,
# corrects a few previously-scraped rows, and appends newly scraped rows.
$wb = $excel.ActiveWorkbook

# --- Sheet "LP1912": refresh scraped schedule rows ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = "Última actualización: 11:19:35"
$ws1.Cells.Item(3,1).Value = "Total filas: 155"

$ws1Rows = @(
  @{ r=6; a="04:03:21"; b="04:46"; c="215A_EL PATO"; d=43; e="LP1912" },
  @{ r=7; a="04:48:48"; b="04:53"; c="11_ETCHEVERRY"; d=5; e="LP1912" },
  @{ r=8; a="04:48:48"; b="05:16"; c="17_ROMERO"; d=28; e="LP1912" },
  @{ r=9; a="04:48:48"; b="05:22"; c="23_HERNANDEZ"; d=34; e="LP1912" },
  @{ r=10; a="05:20:44"; b="05:23"; c="23_HERNANDEZ"; d=3; e="LP1912" },
  @{ r=11; a="05:20:44"; b="05:34"; c="215B_EL PATO"; d=14; e="LP1912" },
  @{ r=12; a="05:20:44"; b="05:46"; c="15_ABASTO"; d=26; e="LP1912" },
  @{ r=13; a="05:20:44"; b="05:54"; c="10_OLMOS"; d=34; e="LP1912" },
  @{ r=14; a="05:20:44"; b="06:04"; c="16_SANTA ANA"; d=44; e="LP1912" },
  @{ r=15; a="06:03:38"; b="06:11"; c="215A_EL PATO"; d=8; e="LP1912" },
  @{ r=16; a="06:03:38"; b="06:14"; c="225_HARAS DEL SUR"; d=11; e="LP1912" },
  @{ r=17; a="06:03:38"; b="06:21"; c="26_HERNANDEZ"; d=18; e="LP1912" },
  @{ r=18; a="06:03:38"; b="06:27"; c="23_HERNANDEZ"; d=24; e="LP1912" },
  @{ r=19; a="06:03:38"; b="06:29"; c="86_EST CHICA-ESC AGRARIA"; d=26; e="LP1912" },
  @{ r=20; a="06:03:38"; b="06:31"; c="16_SANTA ANA"; d=28; e="LP1912" },
  @{ r=21; a="06:03:38"; b="06:44"; c="225_C ROCA-H SUR"; d=41; e="LP1912" },
  @{ r=22; a="06:03:38"; b="06:46"; c="215C_EL PATO"; d=43; e="LP1912" },
  @{ r=23; a="06:52:41"; b="06:56"; c="14_ABASTO"; d=4; e="LP1912" },
  @{ r=24; a="06:03:38"; b="06:59"; c="14_ABASTO"; d=56; e="LP1912" },
  @{ r=25; a="06:52:41"; b="07:01"; c="16_SANTA ANA"; d=9; e="LP1912" },
  @{ r=26; a="06:52:41"; b="07:04"; c="23_HERNANDEZ"; d=12; e="LP1912" },
  @{ r=27; a="06:52:41"; b="07:05"; c="15_ABASTO"; d=13; e="LP1912" },
  @{ r=28; a="06:52:41"; b="07:07"; c="225_GOMEZ"; d=15; e="LP1912" },
  @{ r=29; a="06:52:41"; b="07:11"; c="215A_EL PATO"; d=19; e="LP1912" },
  @{ r=30; a="06:52:41"; b="07:15"; c="11_ETCHEVERRY"; d=23; e="LP1912" },
  @{ r=31; a="06:52:41"; b="07:16"; c="16_SANTA ANA"; d=24; e="LP1912" },
  @{ r=32; a="06:52:41"; b="07:21"; c="26_HERNANDEZ"; d=29; e="LP1912" },
  @{ r=33; a="07:23:38"; b="07:23"; c="16_SANTA ANA"; d=0; e="LP1912" },
  @{ r=34; a="06:52:41"; b="07:23"; c="10_OLMOS"; d=31; e="LP1912" },
  @{ r=35; a="07:23:38"; b="07:31"; c="16_SANTA ANA"; d=8; e="LP1912" },
  @{ r=36; a="07:23:38"; b="07:31"; c="11_ETCHEVERRY"; d=8; e="LP1912" },
  @{ r=37; a="07:23:38"; b="07:32"; c="84_COLONIA URQUIZA-ESC 49"; d=9; e="LP1912" },
  @{ r=38; a="07:23:38"; b="07:36"; c="27_EL RETIRO"; d=13; e="LP1912" },
  @{ r=39; a="07:23:38"; b="07:36"; c="23_HERNANDEZ"; d=13; e="LP1912" },
  @{ r=40; a="07:23:38"; b="07:39"; c="10_OLMOS"; d=16; e="LP1912" },
  @{ r=41; a="07:23:38"; b="07:47"; c="14_ABASTO"; d=24; e="LP1912" },
  @{ r=42; a="07:23:38"; b="07:51"; c="215D_EL PATO"; d=28; e="LP1912" },
  @{ r=43; a="07:23:38"; b="07:55"; c="10_OLMOS"; d=32; e="LP1912" },
  @{ r=44; a="07:59:51"; b="07:59"; c="23_HERNANDEZ"; d=0; e="LP1912" },
  @{ r=45; a="07:59:51"; b="07:59"; c="16_SANTA ANA"; d=0; e="LP1912" },
  @{ r=46; a="07:23:38"; b="08:03"; c="11_ETCHEVERRY"; d=40; e="LP1912" },
  @{ r=47; a="07:59:51"; b="08:04"; c="11_ETCHEVERRY"; d=5; e="LP1912" },
  @{ r=48; a="06:52:41"; b="08:06"; c="23_HERNANDEZ"; d=74; e="LP1912" },
  @{ r=49; a="07:59:51"; b="08:12"; c="15_ABASTO"; d=13; e="LP1912" },
  @{ r=50; a="07:59:51"; b="08:13"; c="10_OLMOS"; d=14; e="LP1912" },
  @{ r=51; a="07:59:51"; b="08:21"; c="26_HERNANDEZ"; d=22; e="LP1912" },
  @{ r=52; a="07:59:51"; b="08:22"; c="16_P MOR-SANTA ANA"; d=23; e="LP1912" },
  @{ r=53; a="07:59:51"; b="08:23"; c="215B_EL PATO"; d=24; e="LP1912" },
  @{ r=54; a="07:59:51"; b="08:27"; c="84_COLONIA URQUIZA-ESC 49"; d=28; e="LP1912" },
  @{ r=55; a="08:31:16"; b="08:32"; c="10_OLMOS"; d=1; e="LP1912" },
  @{ r=56; a="07:59:51"; b="08:34"; c="23_HERNANDEZ"; d=35; e="LP1912" },
  @{ r=57; a="08:31:16"; b="08:35"; c="23_HERNANDEZ"; d=4; e="LP1912" },
  @{ r=58; a="08:31:16"; b="08:42"; c="81_EL PELIGRO"; d=11; e="LP1912" },
  @{ r=59; a="07:59:51"; b="08:43"; c="14_ABASTO"; d=44; e="LP1912" },
  @{ r=60; a="08:31:16"; b="08:44"; c="14_ABASTO"; d=13; e="LP1912" },
  @{ r=61; a="07:23:38"; b="08:44"; c="81_EL PELIGRO"; d=81; e="LP1912" },
  @{ r=62; a="08:31:16"; b="08:53"; c="10_OLMOS"; d=22; e="LP1912" },
  @{ r=63; a="08:54:41"; b="08:54"; c="17_ROMERO"; d=0; e="LP1912" },
  @{ r=64; a="08:54:41"; b="08:57"; c="225_HARAS DEL SUR"; d=3; e="LP1912" },
  @{ r=65; a="08:31:16"; b="09:01"; c="215A_EL PATO"; d=30; e="LP1912" },
  @{ r=66; a="08:54:41"; b="09:03"; c="11_ETCHEVERRY"; d=9; e="LP1912" },
  @{ r=67; a="08:31:16"; b="09:04"; c="11_ETCHEVERRY"; d=33; e="LP1912" },
  @{ r=68; a="08:31:16"; b="09:05"; c="23_HERNANDEZ"; d=34; e="LP1912" },
  @{ r=69; a="08:54:41"; b="09:06"; c="23_HERNANDEZ"; d=12; e="LP1912" },
  @{ r=70; a="08:54:41"; b="09:10"; c="16_P MOR-SANTA ANA"; d=16; e="LP1912" },
  @{ r=71; a="08:31:16"; b="09:11"; c="16_P MOR-SANTA ANA"; d=40; e="LP1912" },
  @{ r=72; a="08:54:41"; b="09:13"; c="10_OLMOS"; d=19; e="LP1912" },
  @{ r=73; a="08:54:41"; b="09:16"; c="27_EL RETIRO"; d=22; e="LP1912" },
  @{ r=74; a="08:31:16"; b="09:17"; c="27_EL RETIRO"; d=46; e="LP1912" },
  @{ r=75; a="07:59:51"; b="09:20"; c="81_EL PELIGRO"; d=81; e="LP1912" },
  @{ r=76; a="08:54:41"; b="09:21"; c="26_HERNANDEZ"; d=27; e="LP1912" },
  @{ r=77; a="07:59:51"; b="09:22"; c="17_ROMERO"; d=83; e="LP1912" },
  @{ r=78; a="08:54:41"; b="09:22"; c="16_SANTA ANA"; d=28; e="LP1912" },
  @{ r=79; a="08:31:16"; b="09:23"; c="16_SANTA ANA"; d=52; e="LP1912" },
  @{ r=80; a="08:54:41"; b="09:23"; c="11_ETCHEVERRY"; d=29; e="LP1912" },
  @{ r=81; a="08:31:16"; b="09:24"; c="11_ETCHEVERRY"; d=53; e="LP1912" },
  @{ r=82; a="08:54:41"; b="09:32"; c="15_ABASTO"; d=38; e="LP1912" },
  @{ r=83; a="09:32:47"; b="09:33"; c="10_OLMOS"; d=1; e="LP1912" },
  @{ r=84; a="09:32:47"; b="09:34"; c="23_HERNANDEZ"; d=2; e="LP1912" },
  @{ r=85; a="08:54:41"; b="09:34"; c="16_SANTA ANA"; d=40; e="LP1912" },
  @{ r=86; a="08:31:16"; b="09:35"; c="16_SANTA ANA"; d=64; e="LP1912" },
  @{ r=87; a="09:32:47"; b="09:41"; c="215C_EL PATO"; d=9; e="LP1912" },
  @{ r=88; a="08:31:16"; b="09:42"; c="215C_EL PATO"; d=71; e="LP1912" },
  @{ r=89; a="09:32:47"; b="09:42"; c="16_SANTA ANA"; d=10; e="LP1912" },
  @{ r=90; a="09:32:47"; b="09:43"; c="14_ABASTO"; d=11; e="LP1912" },
  @{ r=91; a="08:31:16"; b="09:44"; c="14_ABASTO"; d=73; e="LP1912" },
  @{ r=92; a="09:32:47"; b="09:46"; c="16_SANTA ANA"; d=14; e="LP1912" },
  @{ r=93; a="08:54:41"; b="09:52"; c="15_ABASTO"; d=58; e="LP1912" },
  @{ r=94; a="09:32:47"; b="09:53"; c="10_OLMOS"; d=21; e="LP1912" },
  @{ r=95; a="09:32:47"; b="09:54"; c="15_ABASTO"; d=22; e="LP1912" },
  @{ r=96; a="09:32:47"; b="10:03"; c="11_ETCHEVERRY"; d=31; e="LP1912" },
  @{ r=97; a="09:32:47"; b="10:04"; c="23_HERNANDEZ"; d=32; e="LP1912" },
  @{ r=98; a="09:32:47"; b="10:12"; c="15_ABASTO"; d=40; e="LP1912" },
  @{ r=99; a="09:32:47"; b="10:13"; c="10_OLMOS"; d=41; e="LP1912" },
  @{ r=100; a="09:32:47"; b="10:21"; c="26_HERNANDEZ"; d=49; e="LP1912" },
  @{ r=101; a="09:32:47"; b="10:22"; c="17_ROMERO"; d=50; e="LP1912" },
  @{ r=102; a="09:32:47"; b="10:23"; c="11_ETCHEVERRY"; d=51; e="LP1912" },
  @{ r=103; a="09:32:47"; b="10:26"; c="215A_EL PATO"; d=54; e="LP1912" },
  @{ r=104; a="10:39:14"; b="10:40"; c="14_ABASTO"; d=1; e="LP1912" },
  @{ r=105; a="09:32:47"; b="10:41"; c="17_ROMERO"; d=69; e="LP1912" },
  @{ r=106; a="10:39:14"; b="10:42"; c="17_ROMERO"; d=3; e="LP1912" },
  @{ r=107; a="09:32:47"; b="10:43"; c="14_ABASTO"; d=71; e="LP1912" },
  @{ r=108; a="10:39:14"; b="10:47"; c="16_SANTA ANA"; d=8; e="LP1912" },
  @{ r=109; a="10:39:14"; b="10:52"; c="15_ABASTO"; d=13; e="LP1912" },
  @{ r=110; a="10:39:14"; b="10:53"; c="10_OLMOS"; d=14; e="LP1912" },
  @{ r=111; a="09:32:47"; b="10:53"; c="27_EL RETIRO"; d=81; e="LP1912" },
  @{ r=112; a="10:39:14"; b="10:57"; c="27_EL RETIRO"; d=18; e="LP1912" },
  @{ r=113; a="10:39:14"; b="10:57"; c="16_SANTA ANA"; d=18; e="LP1912" },
  @{ r=114; a="09:32:47"; b="11:01"; c="215C_EL PATO"; d=89; e="LP1912" },
  @{ r=115; a="10:39:14"; b="11:02"; c="215C_EL PATO"; d=23; e="LP1912" },
  @{ r=116; a="10:39:14"; b="11:04"; c="11_ETCHEVERRY"; d=25; e="LP1912" },
  @{ r=117; a="10:39:14"; b="11:05"; c="23_HERNANDEZ"; d=26; e="LP1912" },
  @{ r=118; a="09:32:47"; b="11:06"; c="16_P MOR-167 Y 521"; d=94; e="LP1912" },
  @{ r=119; a="10:39:14"; b="11:07"; c="16_P MOR-167 Y 521"; d=28; e="LP1912" },
  @{ r=120; a="10:39:14"; b="11:11"; c="10_OLMOS"; d=32; e="LP1912" },
  @{ r=121; a="10:39:14"; b="11:12"; c="15_ABASTO"; d=33; e="LP1912" },
  @{ r=122; a="11:19:35"; b="11:19"; c="86_EST CHICA-ESC AGRARIA"; d=0; e="LP1912" },
  @{ r=123; a="10:39:14"; b="11:20"; c="86_EST CHICA-ESC AGRARIA"; d=41; e="LP1912" },
  @{ r=124; a="11:19:35"; b="11:21"; c="26_HERNANDEZ"; d=2; e="LP1912" },
  @{ r=125; a="11:19:35"; b="11:22"; c="17_ROMERO"; d=3; e="LP1912" },
  @{ r=126; a="09:32:47"; b="11:26"; c="16_P MOR-SANTA ANA"; d=114; e="LP1912" },
  @{ r=127; a="11:19:35"; b="11:26"; c="16_SANTA ANA"; d=7; e="LP1912" },
  @{ r=128; a="11:19:35"; b="11:27"; c="225_C ROCA-H SUR"; d=8; e="LP1912" },
  @{ r=129; a="11:19:35"; b="11:32"; c="81_EL PELIGRO"; d=13; e="LP1912" },
  @{ r=130; a="11:19:35"; b="11:34"; c="23_HERNANDEZ"; d=15; e="LP1912" },
  @{ r=131; a="11:19:35"; b="11:35"; c="11_ETCHEVERRY"; d=16; e="LP1912" },
  @{ r=132; a="10:39:14"; b="11:35"; c="23_HERNANDEZ"; d=56; e="LP1912" },
  @{ r=133; a="11:19:35"; b="11:36"; c="16_SANTA ANA"; d=17; e="LP1912" },
  @{ r=134; a="10:39:14"; b="11:36"; c="11_ETCHEVERRY"; d=57; e="LP1912" },
  @{ r=135; a="11:19:35"; b="11:42"; c="17_ROMERO"; d=23; e="LP1912" },
  @{ r=136; a="11:19:35"; b="11:43"; c="10_OLMOS"; d=24; e="LP1912" },
  @{ r=137; a="11:19:35"; b="11:51"; c="215B_EL PATO"; d=32; e="LP1912" },
  @{ r=138; a="11:19:35"; b="11:52"; c="15_ABASTO"; d=33; e="LP1912" },
  @{ r=139; a="11:19:35"; b="11:59"; c="225_GOMEZ"; d=40; e="LP1912" },
  @{ r=140; a="10:39:14"; b="12:02"; c="84_COLONIA URQUIZA-ESC 49"; d=83; e="LP1912" },
  @{ r=141; a="11:19:35"; b="12:04"; c="23_HERNANDEZ"; d=45; e="LP1912" },
  @{ r=142; a="11:19:35"; b="12:06"; c="84_COLONIA URQUIZA-ESC 49"; d=47; e="LP1912" },
  @{ r=143; a="11:19:35"; b="12:06"; c="14_ABASTO"; d=47; e="LP1912" },
  @{ r=144; a="11:19:35"; b="12:06"; c="16_P MOR-SANTA ANA"; d=47; e="LP1912" },
  @{ r=145; a="10:39:14"; b="12:07"; c="16_P MOR-SANTA ANA"; d=88; e="LP1912" },
  @{ r=146; a="11:19:35"; b="12:13"; c="10_OLMOS"; d=54; e="LP1912" },
  @{ r=147; a="11:19:35"; b="12:14"; c="17_ROMERO"; d=55; e="LP1912" },
  @{ r=148; a="11:19:35"; b="12:20"; c="215A_EL PATO"; d=61; e="LP1912" },
  @{ r=149; a="11:19:35"; b="12:20"; c="14_ABASTO"; d=61; e="LP1912" },
  @{ r=150; a="10:39:14"; b="12:21"; c="215A_EL PATO"; d=102; e="LP1912" },
  @{ r=151; a="11:19:35"; b="12:21"; c="26_HERNANDEZ"; d=62; e="LP1912" },
  @{ r=152; a="10:39:14"; b="12:22"; c="14_ABASTO"; d=103; e="LP1912" },
  @{ r=153; a="11:19:35"; b="12:36"; c="27_EL RETIRO"; d=77; e="LP1912" },
  @{ r=154; a="10:39:14"; b="12:37"; c="27_EL RETIRO"; d=118; e="LP1912" },
  @{ r=155; a="11:19:35"; b="12:38"; c="17_179 Y 38"; d=79; e="LP1912" },
  @{ r=156; a="11:19:35"; b="12:41"; c="10_OLMOS"; d=82; e="LP1912" },
  @{ r=157; a="11:19:35"; b="12:48"; c="11_ETCHEVERRY"; d=89; e="LP1912" },
  @{ r=158; a="11:19:35"; b="12:50"; c="15_ABASTO"; d=91; e="LP1912" },
  @{ r=159; a="11:19:35"; b="13:06"; c="16_P MOR-SANTA ANA"; d=107; e="LP1912" },
  @{ r=160; a="11:19:35"; b="13:14"; c="215D_EL PATO"; d=115; e="LP1912" }
)

foreach ($row in $ws1Rows) {
  $ws1.Cells.Item($row.r, 1).Value = $row.a
  $ws1.Cells.Item($row.r, 2).Value = $row.b
  $ws1.Cells.Item($row.r, 3).Value = $row.c
  $ws1.Cells.Item($row.r, 4).Value = $row.d
  $ws1.Cells.Item($row.r, 5).Value = $row.e
}

# --- Sheet "LP1912-215": refresh scraped schedule rows ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 11:19:35"
$ws2.Cells.Item(3,1).Value = "Total filas: 17"

$ws2Rows = @(
  @{ r=6; a="04:03:21"; b="04:46"; c="215A_EL PATO"; d=43; e="LP1912" },
  @{ r=7; a="05:20:44"; b="05:34"; c="215B_EL PATO"; d=14; e="LP1912" },
  @{ r=8; a="06:03:38"; b="06:11"; c="215A_EL PATO"; d=8; e="LP1912" },
  @{ r=9; a="06:03:38"; b="06:46"; c="215C_EL PATO"; d=43; e="LP1912" },
  @{ r=10; a="06:52:41"; b="07:11"; c="215A_EL PATO"; d=19; e="LP1912" },
  @{ r=11; a="07:23:38"; b="07:51"; c="215D_EL PATO"; d=28; e="LP1912" },
  @{ r=12; a="07:59:51"; b="08:23"; c="215B_EL PATO"; d=24; e="LP1912" },
  @{ r=13; a="08:31:16"; b="09:01"; c="215A_EL PATO"; d=30; e="LP1912" },
  @{ r=14; a="09:32:47"; b="09:41"; c="215C_EL PATO"; d=9; e="LP1912" },
  @{ r=15; a="08:31:16"; b="09:42"; c="215C_EL PATO"; d=71; e="LP1912" },
  @{ r=16; a="09:32:47"; b="10:26"; c="215A_EL PATO"; d=54; e="LP1912" },
  @{ r=17; a="09:32:47"; b="11:01"; c="215C_EL PATO"; d=89; e="LP1912" },
  @{ r=18; a="10:39:14"; b="11:02"; c="215C_EL PATO"; d=23; e="LP1912" },
  @{ r=19; a="11:19:35"; b="11:51"; c="215B_EL PATO"; d=32; e="LP1912" },
  @{ r=20; a="11:19:35"; b="12:20"; c="215A_EL PATO"; d=61; e="LP1912" },
  @{ r=21; a="10:39:14"; b="12:21"; c="215A_EL PATO"; d=102; e="LP1912" },
  @{ r=22; a="11:19:35"; b="13:14"; c="215D_EL PATO"; d=115; e="LP1912" }
)

foreach ($row in $ws2Rows) {
  $ws2.Cells.Item($row.r, 1).Value = $row.a
  $ws2.Cells.Item($row.r, 2).Value = $row.b
  $ws2.Cells.Item($row.r, 3).Value = $row.c
  $ws2.Cells.Item($row.r, 4).Value = $row.d
  $ws2.Cells.Item($row.r, 5).Value = $row.e
}

# --- Sheet "6203-6173": refresh scraped schedule rows ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: 11:19:35"
$ws3.Cells.Item(3,1).Value = "Total filas: 22"

$ws3Rows = @(
  @{ r=6; a="05:20:44"; b="05:43"; c="215A_LA PLATA"; d=23; e="L6173" },
  @{ r=7; a="04:03:21"; b="05:44"; c="215A_LA PLATA"; d=101; e="L6173" },
  @{ r=8; a="06:03:38"; b="06:08"; c="215A_LA PLATA"; d=5; e="L6173" },
  @{ r=9; a="06:03:38"; b="06:32"; c="215C_LA PLATA"; d=29; e="L6203" },
  @{ r=10; a="06:52:41"; b="06:59"; c="215B_LP-P MOR-1 Y 57"; d=7; e="L6173" },
  @{ r=11; a="06:03:38"; b="07:00"; c="215B_LP-P MOR-1 Y 57"; d=57; e="L6173" },
  @{ r=12; a="07:23:38"; b="07:35"; c="215A_LA PLATA"; d=12; e="L6173" },
  @{ r=13; a="06:52:41"; b="08:06"; c="215C_LA PLATA"; d=74; e="L6203" },
  @{ r=14; a="07:23:38"; b="08:09"; c="215C_LA PLATA"; d=46; e="L6203" },
  @{ r=15; a="07:59:51"; b="08:11"; c="215C_LA PLATA"; d=12; e="L6203" },
  @{ r=16; a="06:52:41"; b="08:31"; c="215A_LA PLATA"; d=99; e="L6173" },
  @{ r=17; a="07:59:51"; b="08:35"; c="215A_LA PLATA"; d=36; e="L6173" },
  @{ r=18; a="08:31:16"; b="08:37"; c="215A_LA PLATA"; d=6; e="L6173" },
  @{ r=19; a="07:59:51"; b="09:08"; c="215D_LA PLATA"; d=69; e="L6203" },
  @{ r=20; a="08:54:41"; b="09:09"; c="215D_LA PLATA"; d=15; e="L6203" },
  @{ r=21; a="09:32:47"; b="10:02"; c="215B_LP-P MOR-40 Y 115"; d=30; e="L6173" },
  @{ r=22; a="08:31:16"; b="10:03"; c="215B_LP-P MOR-40 Y 115"; d=92; e="L6173" },
  @{ r=23; a="10:39:14"; b="10:54"; c="215A_LA PLATA"; d=15; e="L6173" },
  @{ r=24; a="09:32:47"; b="11:13"; c="215C_LA PLATA"; d=101; e="L6203" },
  @{ r=25; a="10:39:14"; b="11:14"; c="215C_LA PLATA"; d=35; e="L6203" },
  @{ r=26; a="11:19:35"; b="12:04"; c="215A_LA PLATA"; d=45; e="L6173" },
  @{ r=27; a="11:19:35"; b="12:53"; c="215C_LA PLATA"; d=94; e="L6203" }
)

foreach ($row in $ws3Rows) {
  $ws3.Cells.Item($row.r, 1).Value = $row.a
  $ws3.Cells.Item($row.r, 2).Value = $row.b
  $ws3.Cells.Item($row.r, 3).Value = $row.c
  $ws3.Cells.Item($row.r, 4).Value = $row.d
  $ws3.Cells.Item($row.r, 5).Value = $row.e
}

